$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.010.39"
$ws.Range("E2").Value = '  -1.82%  '

$ws.Range("D3").Value = "'1.555.06"
$ws.Range("E3").Value = '  -0.82%  '

$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("D6").Value = "'286.99"
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("D7").Value = "'0.3817"
$ws.Range("E7").Value = '  +3.40%  '

$ws.Range("D8").Value = "'0.3233"
$ws.Range("E8").Value = '  -1.87%  '

$ws.Range("D9").Value = "'41.35"
$ws.Range("E9").Value = '  -12.35%  '

$ws.Range("E10").Value = '  -2.89%  '

$ws.Range("D11").Value = "'0.07311"
$ws.Range("E11").Value = '  -1.80%  '

$ws.Range("E12").Value = '  -0.13%  '

$ws.Range("D13").Value = "'19.36"
$ws.Range("E13").Value = '  -6.29%  '

$ws.Range("D14").Value = "'5.715"
$ws.Range("E14").Value = '  -3.21%  '

$ws.Range("D15").Value = "'6.813"
$ws.Range("E15").Value = '  -0.51%  '

$ws.Range("D16").Value = "'1.555.94"
$ws.Range("E16").Value = '  -0.53%  '

$ws.Range("E17").Value = '  -1.36%  '

$ws.Range("D18").Value = "'0.06625"
$ws.Range("E18").Value = '  -1.22%  '

$ws.Range("D19").Value = "'85.22"
$ws.Range("E19").Value = '  -1.83%  '

$ws.Range("D20").Value = "'6.406"
$ws.Range("E20").Value = '  +0.75%  '

$ws.Range("D21").Value = "'0.9998"
$ws.Range("E21").Value = '  -0.10%  '

$ws.Range("E22").Value = '  -2.75%  '

$ws.Range("D23").Value = "'11.43"
$ws.Range("E23").Value = '  -3.87%  '

$ws.Range("D24").Value = "'22.009.65"
$ws.Range("E24").Value = '  -1.74%  '

$ws.Range("D25").Value = "'2.292"
$ws.Range("E25").Value = '  -2.93%  '

$ws.Range("D26").Value = "'2.522"
$ws.Range("E26").Value = '  -2.87%  '

$ws.Range("D27").Value = "'148.78"
$ws.Range("E27").Value = '  -1.38%  '

$ws.Range("D28").Value = "'18.81"
$ws.Range("E28").Value = '  -3.32%  '

$ws.Range("D29").Value = "'4.857"
$ws.Range("E29").Value = '  -1.59%  '

$ws.Range("D30").Value = "'1.727.54"
$ws.Range("E30").Value = '  -1.06%  '

$ws.Range("D31").Value = "'120.54"
$ws.Range("E31").Value = '  -2.93%  '

$ws.Range("D32").Value = "'1.092"
$ws.Range("E32").Value = '  +1.71%  '

$ws.Range("D33").Value = "'5.872"
$ws.Range("E33").Value = '  -2.57%  '

$ws.Range("D34").Value = "'9.268"
$ws.Range("E34").Value = '  -5.42%  '

$ws.Range("D35").Value = "'0.08127"
$ws.Range("E35").Value = '  -1.99%  '

$ws.Range("D36").Value = "'1.642"
$ws.Range("E36").Value = '  -16.84%  '

$ws.Range("D37").Value = "'0.06206"
$ws.Range("E37").Value = '  -2.32%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = "'0.02296"
$ws.Range("E38").Value = '  -5.01%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = "'5.244"
$ws.Range("E39").Value = '  -0.10%  '

$ws.Range("D40").Value = "'0.2104"
$ws.Range("E40").Value = '  -4.16%  '

$ws.Range("D41").Value = "'1.219"
$ws.Range("E41").Value = '  -5.51%  '

$ws.Range("D42").Value = "'10.84"
$ws.Range("E42").Value = '  -3.97%  '

$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("D44").Value = "'0.5924"
$ws.Range("E44").Value = '  -3.45%  '

$ws.Range("D45").Value = "'13.47"
$ws.Range("E45").Value = '  -2.91%  '

$ws.Range("D47").Value = "'0.5737"
$ws.Range("E47").Value = '  -4.07%  '

$ws.Range("D48").Value = "'1.930"
$ws.Range("E48").Value = '  -4.70%  '

$ws.Range("D49").Value = "'119.33"
$ws.Range("E49").Value = '  -4.28%  '

$ws.Range("D50").Value = "'1.154"
$ws.Range("E50").Value = '  -3.26%  '

$ws.Range("D51").Value = "'0.06867"
$ws.Range("E51").Value = '  -4.20%  '
